# Applies the crypto price/volume/listing update described by the commit:
# "Updated cryptos list on Sat Oct 19 11:58:19 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.125.52"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.635.44"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D5").Value = "'597.80"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "'154.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.543"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "2.633.82"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "'5.23"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "'0.350"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "'27.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "3.114.51"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "67.997.57"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "2.662.14"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "'11.32"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'362.27"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").Value = "'7.38"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'4.34"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("D23").Value = "'4.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "'75.08"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.40%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'9.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "'561.25"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").Value = "'7.99"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").Value = "'0.129"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'1.58"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.86%  "
$ws.Range("D38").Value = "'160.78"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").Value = "'19.23"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("D41").Value = "'1.87"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'5.29"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "0.0₆0338"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D47").Value = "'40.56"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").Value = "'156.81"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0785"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'21.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.36%  "

Write-Host "Applied 85 cell updates"
